# Add US01 & US02 test-mapping details to the Sprint1 sheet's sprint report
# (columns G-O for rows 2 and 3), mirroring the existing US07/US08 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# Row 2 -> US01 "Dates before current date"
$ws.Range("G2").Value = 27
$ws.Range("H2").Value = 30
$ws.Range("I2").Value = "Yes"
$ws.Range("K2").Value = "US01()"
$ws.Range("L2").Value = "371-401"
$ws.Range("N2").Value = "test_userstory01()"
$ws.Range("O2").Value = "18-20"

# Row 3 -> US02 "Birth before marriage"
$ws.Range("G3").Value = 13
$ws.Range("H3").Value = 20
$ws.Range("I3").Value = "Yes"
$ws.Range("K3").Value = "US02()"
$ws.Range("L3").Value = "403 - 416"
$ws.Range("N3").Value = "test_userstory02()"
$ws.Range("O3").Value = "22-24"

# Update sheet view to match the author's new selection/scroll position
$ws.Activate()
$ws.Range("O3").Select()
$excel.ActiveWindow.ScrollColumn = 5
